# Update Common IC Type counts for 2010 data (column B, rows 2-32),
# label + recount the former "blank" bucket in row 33, and append a new
# "unpopulated" bucket row 34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3394
    3  = 420
    4  = 36
    5  = 70
    6  = 722
    7  = 100
    8  = 180
    9  = 2030
    10 = 96
    11 = 9344
    12 = 310
    13 = 190
    14 = 10332
    15 = 722
    16 = 546
    17 = 1144
    18 = 150
    19 = 670
    20 = 698
    21 = 2742
    22 = 4836
    23 = 438
    24 = 10
    25 = 10
    26 = 70
    27 = 198
    28 = 166
    29 = 950
    30 = 288
    31 = 4896
    32 = 14326
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# Row 33 previously held an unlabeled count; it now gets the "unpopulated"
# label and a new, much smaller count.
$ws.Cells.Item(33, 1).Value = "unpopulated"
$ws.Cells.Item(33, 2).Value = 6

# New row 34 holds the remaining unlabeled count.
$ws.Cells.Item(34, 1).Value = ""
$ws.Cells.Item(34, 2).Value = 706
